$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'245.73"

# Row 3
$ws.Range("D3").Value = "'23.94"

# Row 4
$ws.Range("D4").Value = "'5.153"

# Row 5
$ws.Range("D5").Value = "'0.05740"

# Row 6
$ws.Range("D6").Value = "'6.486"

# Row 7
$ws.Range("D7").Value = "'3.167"

# Row 9
$ws.Range("D9").Value = "'0.8538"

# Row 10
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = "'0.1379"
$ws.Range("E10").Value = '9WazirXWRX'

# Row 11
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = "'0.06946"
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'

# Row 12
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").Value = "'0.03179"
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'

# Row 13
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = "'0.02883"
$ws.Range("E13").Value = '12BitrueCoinBTR'

# Row 14
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = "'0.09350"
$ws.Range("E14").Value = '13BitMartTokenBMX'

# Row 15
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").Value = "'3.821"
$ws.Range("E15").Value = '14MCDexMCB'

# Row 16
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").Value = "'0.001521"
$ws.Range("E16").Value = '15BitForexTokenBF'

# Row 17
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").Value = "'0.04697"
$ws.Range("E17").Value = '16CoinExTokenCET'

# Row 18
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").Value = "'0.0005969"
$ws.Range("E18").Value = '17OneONE'

# Row 19
$ws.Range("D19").Value = "'0.006187"

# Row 20
$ws.Range("D20").Value = "'0.001242"

# Row 21
$ws.Range("D21").Value = "'0.004825"

# Row 22
$ws.Range("B22").Value = 'UpBots'
$ws.Range("C22").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D22").Value = "'0.007488"
$ws.Range("E22").Value = '21UpBotsUBXTBestin24h'

# Row 23
$ws.Range("B23").Value = 'NitroEx'
$ws.Range("C23").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D23").Value = "'0.00008495"
$ws.Range("E23").Value = '22NitroExNTX'

# Row 24
$ws.Range("B24").Value = 'LEO'
$ws.Range("C24").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D24").Value = "'3.541"
$ws.Range("E24").Value = '23LEOLEO'

# Row 25
$ws.Range("B25").Value = 'BTSEToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D25").Value = "'2.153"
$ws.Range("E25").Value = '24BTSETokenBTSE'

# Row 26
$ws.Range("B26").Value = 'BitpandaEcosystemToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D26").Value = "'0.3199"
$ws.Range("E26").Value = '25BitpandaEcosystemTokenBEST'

# Row 27
$ws.Range("B27").Value = 'ProBitToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D27").Value = "'0.1338"
$ws.Range("E27").Value = '26ProBitTokenPROB'

# Row 40
$ws.Range("D40").Value = "'0.03693"

# Row 41
$ws.Range("D41").Value = "'0.006366"

# Row 42
$ws.Range("D42").Value = "'0.1056"

# Row 43
$ws.Range("D43").Value = "'0.002259"

# Row 44
$ws.Range("D44").Value = "'0.007817"

# Row 45
$ws.Range("D45").Value = "'0.00005461"

# Row 46
$ws.Range("D46").Value = "'0.00000000750"

# Row 47
$ws.Range("D47").Value = "'0.4000"

# Row 48
$ws.Range("D48").Value = "'0.002450"

# Row 49
$ws.Range("D49").Value = "'0.00002100"

# Row 50
$ws.Range("D50").Value = "'0.0002000"
